$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new row 2 for the
#    2022-Q3 quarter and renumber the index column / shift the rest
#    down (Excel does the shifting for us via Rows.Insert()).
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# New row's data cells shouldn't inherit the blank-row style that
# Insert() stamped on them - clear it so they match the plain data
# cells below (no explicit style), then copy the index-column style
# from the row underneath onto the new A2 cell.
$summary.Range("B2:D2").ClearFormats()
$summary.Range("A3").Copy($summary.Range("A2"))

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 6
$summary.Range("D2").Value = 0.54

# The rows that got pushed down keep their old index-column values,
# so renumber A3:A8 sequentially (1..6) to match the re-flowed order.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6

# ---------------------------------------------------------------
# 2. Create the new "2022-Q3" sheet carrying the fund holdings for
#    that quarter. Copy an existing quarter sheet with the same
#    shape (header + formatting) so the new sheet inherits identical
#    styling, then trim/overwrite its contents. "2022-Q1" has 8 data
#    rows; trim the last 2 to get down to the 6 rows we need, then
#    place the copy immediately before "2022-Q1" (i.e. sheet index 2).
# ---------------------------------------------------------------
$sourceSheet = $wb.Worksheets.Item("2022-Q1")
$beforeSheet = $wb.Worksheets.Item("2022-Q1")
$sourceSheet.Copy($beforeSheet, $null)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Drop the 2 extra data rows (rows 8 and 9) so only 6 remain (rows 2-7).
$q3.Rows.Item(9).Delete()
$q3.Rows.Item(8).Delete()

# --- Row 2 : 016935 / 景顺长城中证500指数增强C ---
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "'016935"
$q3.Range("C2").Value = "景顺长城中证500指数增强C"
$q3.Range("D2").Value = "'15.57"
$q3.Range("E2").Value = "'93.89"
$q3.Range("F2").Value = "'2.02"
$q3.Range("G2").Value = "'0.3145"
$q3.Range("H2").Value = 7

# --- Row 3 : 005994 / 国投瑞银中证500指数量化增强A ---
$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "'005994"
$q3.Range("C3").Value = "国投瑞银中证500指数量化增强A"
$q3.Range("D3").Value = "'13.36"
$q3.Range("E3").Value = "'88.67"
$q3.Range("F3").Value = "'1.22"
$q3.Range("G3").Value = "'0.1630"
$q3.Range("H3").Value = 4

# --- Row 4 : 007089 / 国投瑞银中证500指数量化增强C ---
$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "'007089"
$q3.Range("C4").Value = "国投瑞银中证500指数量化增强C"
$q3.Range("D4").Value = "'4.45"
$q3.Range("E4").Value = "'88.67"
$q3.Range("F4").Value = "'1.22"
$q3.Range("G4").Value = "'0.0543"
$q3.Range("H4").Value = 4

# --- Row 5 : 161038 / 富国新兴成长量化精选混合（LOF）A ---
$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "'161038"
$q3.Range("C5").Value = "富国新兴成长量化精选混合（LOF）A"
$q3.Range("D5").Value = "'0.81"
$q3.Range("E5").Value = "'91.98"
$q3.Range("F5").Value = "'1.37"
$q3.Range("G5").Value = "'0.0111"
$q3.Range("H5").Value = 6

# --- Row 6 : 014171 / 富国新兴成长量化精选混合（LOF）C ---
$q3.Range("A6").Value = 4
$q3.Range("B6").Value = "'014171"
$q3.Range("C6").Value = "富国新兴成长量化精选混合（LOF）C"
$q3.Range("D6").Value = "'0.00"
$q3.Range("E6").Value = "'91.98"
$q3.Range("F6").Value = "'1.37"
$q3.Range("G6").Value = 0
$q3.Range("H6").Value = 6

# --- Row 7 : 006682 / 景顺长城中证500指数增强A ---
$q3.Range("A7").Value = 5
$q3.Range("B7").Value = "'006682"
$q3.Range("C7").Value = "景顺长城中证500指数增强A"
$q3.Range("D7").Value = "'0.00"
$q3.Range("E7").Value = "'93.89"
$q3.Range("F7").Value = "'2.02"
$q3.Range("G7").Value = 0
$q3.Range("H7").Value = 7

Write-Host "2022-Q3 sheet created and 总计 sheet updated"
